$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Range adjustment percentage" label to mention the effective date.
$ws.Range("A3").Value = "Range adjustment percentage eff 7/1"

# Add a new row documenting the department approver.
$ws.Range("A9").Value = "Dept Approver"
$ws.Range("B9").Value = "Ms. Norbury, Mathletes Coach"

# Widen column A so the longer labels are no longer truncated.
$ws.Columns.Item(1).ColumnWidth = 30.8

# Update the sheet's stored selection.
[void]$ws.Range("A12").Select()
